# Slide 2, "Content Placeholder 2": split the last bullet's single run into
# three runs so the middle portion ("implemented at" -> "implemented in")
# becomes its own run, matching the target OOXML:
#   "Communication is " | "implemented in " | "the application layer using  a set of well-known LWM2M types"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Locate the target paragraph by its current (pre-edit) text so the script
# is resilient to any re-ordering of the other bullets.
$targetOld = "Communication is implemented at the application layer using  a set of well-known LWM2M types"
$paraCount = $tr.Paragraphs().Count
$paraIndex = -1
for ($i = 1; $i -le $paraCount; $i++) {
    if ($tr.Paragraphs($i).Text -eq $targetOld) {
        $paraIndex = $i
        break
    }
}

$para = $tr.Paragraphs($paraIndex)
$fullText = $para.Text

$oldMiddle = "implemented at "
$newMiddle = "implemented in "

$startZero = $fullText.IndexOf($oldMiddle)
$startOne  = $startZero + 1
$len       = $oldMiddle.Length

# Re-point the middle slice's text; this splits the original single run
# into three runs (before / middle / after) while leaving the untouched
# leading and trailing text runs' formatting intact.
$midRange = $para.Characters($startOne, $len)
$midRange.Text = $newMiddle
